$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New OHLCV rows appended after the existing last row (1084), extending
# the data through 2025-10-0x timestamps (rows 1085:1116).
$numRows = 32
$numCols = 6
$data = New-Object 'object[,]' $numRows,$numCols

$data[0,0] = 45534.5
$data[0,1] = 11.03
$data[0,2] = 11.17
$data[0,3] = 10.51
$data[0,4] = 10.54
$data[0,5] = 550427.75
$data[1,0] = 45534.66666666666
$data[1,1] = 10.54
$data[1,2] = 11.11
$data[1,3] = 10.5
$data[1,4] = 11
$data[1,5] = 367917.74
$data[2,0] = 45534.83333333334
$data[2,1] = 10.99
$data[2,2] = 11.14
$data[2,3] = 10.96
$data[2,4] = 11.11
$data[2,5] = 163267.88
$data[3,0] = 45535
$data[3,1] = 11.12
$data[3,2] = 11.18
$data[3,3] = 11.07
$data[3,4] = 11.11
$data[3,5] = 125490.37
$data[4,0] = 45535.16666666666
$data[4,1] = 11.12
$data[4,2] = 11.18
$data[4,3] = 11.09
$data[4,4] = 11.12
$data[4,5] = 87331.38
$data[5,0] = 45535.33333333334
$data[5,1] = 11.13
$data[5,2] = 11.2
$data[5,3] = 11.03
$data[5,4] = 11.15
$data[5,5] = 159002.15
$data[6,0] = 45535.5
$data[6,1] = 11.14
$data[6,2] = 11.21
$data[6,3] = 11.03
$data[6,4] = 11.09
$data[6,5] = 131788.48
$data[7,0] = 45535.66666666666
$data[7,1] = 11.09
$data[7,2] = 11.11
$data[7,3] = 10.93
$data[7,4] = 11.03
$data[7,5] = 177944.94
$data[8,0] = 45535.83333333334
$data[8,1] = 11.02
$data[8,2] = 11.06
$data[8,3] = 10.99
$data[8,4] = 11.06
$data[8,5] = 28981.75
$data[9,0] = 45536
$data[9,1] = 11.02
$data[9,2] = 11.02
$data[9,3] = 10.83
$data[9,4] = 10.94
$data[9,5] = 165212.69
$data[10,0] = 45536.16666666666
$data[10,1] = 10.94
$data[10,2] = 10.97
$data[10,3] = 10.78
$data[10,4] = 10.88
$data[10,5] = 186669.68
$data[11,0] = 45536.33333333334
$data[11,1] = 10.88
$data[11,2] = 10.9
$data[11,3] = 10.7
$data[11,4] = 10.81
$data[11,5] = 165189.02
$data[12,0] = 45536.5
$data[12,1] = 10.8
$data[12,2] = 10.82
$data[12,3] = 10.52
$data[12,4] = 10.7
$data[12,5] = 348491.28
$data[13,0] = 45536.66666666666
$data[13,1] = 10.7
$data[13,2] = 10.74
$data[13,3] = 10.69
$data[13,4] = 10.7
$data[13,5] = 16914.8
$data[14,0] = 45536.83333333334
$data[14,1] = 10.71
$data[14,2] = 10.77
$data[14,3] = 10.23
$data[14,4] = 10.37
$data[14,5] = 386833.01
$data[15,0] = 45537
$data[15,1] = 10.37
$data[15,2] = 10.51
$data[15,3] = 10.32
$data[15,4] = 10.48
$data[15,5] = 190812.54
$data[16,0] = 45537.16666666666
$data[16,1] = 10.48
$data[16,2] = 10.52
$data[16,3] = 10.33
$data[16,4] = 10.38
$data[16,5] = 171340.48
$data[17,0] = 45537.33333333334
$data[17,1] = 10.38
$data[17,2] = 10.79
$data[17,3] = 10.37
$data[17,4] = 10.74
$data[17,5] = 295421.94
$data[18,0] = 45537.5
$data[18,1] = 10.75
$data[18,2] = 10.79
$data[18,3] = 10.59
$data[18,4] = 10.66
$data[18,5] = 308048.74
$data[19,0] = 45537.66666666666
$data[19,1] = 10.66
$data[19,2] = 10.71
$data[19,3] = 10.55
$data[19,4] = 10.66
$data[19,5] = 190022.86
$data[20,0] = 45537.83333333334
$data[20,1] = 10.67
$data[20,2] = 10.84
$data[20,3] = 10.64
$data[20,4] = 10.77
$data[20,5] = 159313.35
$data[21,0] = 45538
$data[21,1] = 10.77
$data[21,2] = 10.88
$data[21,3] = 10.69
$data[21,4] = 10.7
$data[21,5] = 158383.22
$data[22,0] = 45538.16666666666
$data[22,1] = 10.7
$data[22,2] = 10.73
$data[22,3] = 10.61
$data[22,4] = 10.69
$data[22,5] = 174493.44
$data[23,0] = 45538.33333333334
$data[23,1] = 10.7
$data[23,2] = 10.71
$data[23,3] = 10.54
$data[23,4] = 10.67
$data[23,5] = 186015.58
$data[24,0] = 45538.5
$data[24,1] = 10.67
$data[24,2] = 10.74
$data[24,3] = 10.27
$data[24,4] = 10.3
$data[24,5] = 439936.79
$data[25,0] = 45538.66666666666
$data[25,1] = 10.29
$data[25,2] = 10.51
$data[25,3] = 10.21
$data[25,4] = 10.5
$data[25,5] = 389523.58
$data[26,0] = 45538.83333333334
$data[26,1] = 10.49
$data[26,2] = 10.52
$data[26,3] = 10.29
$data[26,4] = 10.31
$data[26,5] = 248381.63
$data[27,0] = 45539
$data[27,1] = 10.32
$data[27,2] = 10.39
$data[27,3] = 9.82
$data[27,4] = 10.25
$data[27,5] = 697334.65
$data[28,0] = 45539.16666666666
$data[28,1] = 10.25
$data[28,2] = 10.41
$data[28,3] = 10.17
$data[28,4] = 10.26
$data[28,5] = 312608
$data[29,0] = 45539.33333333334
$data[29,1] = 10.27
$data[29,2] = 10.3
$data[29,3] = 10.09
$data[29,4] = 10.12
$data[29,5] = 241012.14
$data[30,0] = 45539.5
$data[30,1] = 10.12
$data[30,2] = 10.41
$data[30,3] = 10.06
$data[30,4] = 10.4
$data[30,5] = 561812.48
$data[31,0] = 45539.66666666666
$data[31,1] = 10.39
$data[31,2] = 10.5
$data[31,3] = 10.27
$data[31,4] = 10.3
$data[31,5] = 316415.81

$startRow = 1085
$endRow = $startRow + $numRows - 1
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, $numCols))
$rng.Value = $data

# Column A uses a custom date/time style (border + bold + centered + the
# "YYYY-MM-DD HH:MM:SS" number format) throughout the sheet. Copy that exact
# cell style from the previous last row onto the new column-A cells so the
# appended rows look identical to the rest of the column.
$srcA = $ws.Cells.Item(1084, 1)
$dstA = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))
$srcA.Copy()
$dstA.PasteSpecial(-4122)
$excel.CutCopyMode = $false
